$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Hendrik De Wilde"

$ws.Range("B12").Value = 44623
$ws.Range("B12").NumberFormat = "d-mmm"

$ws.Range("C12").Value = 0.75
$ws.Range("C12").NumberFormat = "h:mm"

$ws.Range("D12").Value = 0.83333333333333337
$ws.Range("D12").NumberFormat = "h:mm"

$ws.Range("E12").Value = 2

$ws.Range("F12").Value = "Validation of NationRegistrationNumber en Constructor Models"

$ws.Range("F12").Select()
